$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing header cell (AC1) onto the
# three new header cells so they match the bold/centered/bordered style
# used by the rest of row 1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Header row (row 1) - new columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Data rows 2-46: team record values, same for every player row
for ($r = 2; $r -le 46; $r++) {
    $ws.Cells.Item($r, 30).Value = 82
    $ws.Cells.Item($r, 31).Value = 80
    $ws.Cells.Item($r, 32).Value = 0
}
